$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.275.88'
$ws.Range("E2").Value = '  +1.48%  '

$ws.Range("D3").Value = '3.875.11'
$ws.Range("E3").Value = '  +0.92%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").Value = '''471.30'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +10.19%  '

$ws.Range("D6").Value = '''145.43'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +10.34%  '

$ws.Range("D7").Value = '''0.636'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +3.32%  '

$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("D9").Value = '''0.749'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.13%  '

$ws.Range("D10").Value = '''0.154'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.05%  '

$ws.Range("D11").Value = '''0.0000312'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -8.19%  '

$ws.Range("D12").Value = '''43.67'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +4.01%  '

$ws.Range("D13").Value = '''10.45'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.00%  '

$ws.Range("D14").Value = '4.496.12'
$ws.Range("E14").Value = '  +0.90%  '

$ws.Range("D15").Value = '''14.82'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -5.94%  '

$ws.Range("D16").Value = '3.868.37'
$ws.Range("E16").Value = '  +0.70%  '

$ws.Range("E17").Value = '  -0.32%  '

$ws.Range("D18").Value = '''20.09'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.09%  '

$ws.Range("D19").Value = '''1.17'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +6.08%  '

$ws.Range("D20").Value = '67.449.28'
$ws.Range("E20").Value = '  +1.27%  '

$ws.Range("D21").Value = '''436.43'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +4.85%  '

$ws.Range("D22").Value = '''14.94'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.85%  '

$ws.Range("D23").Value = '''3.31'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +6.06%  '

$ws.Range("D24").Value = '''89.35'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.95%  '

$ws.Range("D25").Value = '''3.61'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +9.66%  '

$ws.Range("D26").Value = '''38.06'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.54%  '

$ws.Range("D27").Value = '''10.19'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +7.90%  '

$ws.Range("D28").Value = '''10.00'
$ws.Range("D28").ClearFormats()

$ws.Range("D29").Value = '''5.48'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.11%  '

$ws.Range("D30").Value = '''734.15'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.57%  '

$ws.Range("D31").Value = '''13.90'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.19%  '

$ws.Range("D32").Value = '''0.135'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +6.87%  '

$ws.Range("E33").Value = '  +2.93%  '

$ws.Range("D34").Value = '''44.71'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +13.68%  '

$ws.Range("D35").Value = '''0.164'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +8.22%  '

$ws.Range("D36").Value = '''58.09'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.35%  '

$ws.Range("D37").Value = '''1.00'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.10%  '

$ws.Range("D38").Value = '''5.50'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.43%  '

$ws.Range("D39").Value = '''0.0485'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.52%  '

$ws.Range("E40").Value = '  +8.36%  '

$ws.Range("D41").Value = '''2.91'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.33%  '

$ws.Range("D42").Value = '0.0₃0691'
$ws.Range("E42").Value = '  -6.99%  '

$ws.Range("E43").Value = '  +3.20%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '''1.00'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.06%  '

$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").Value = '''2.56'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +7.13%  '

$ws.Range("D46").Value = '''3.47'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.89%  '

$ws.Range("D47").Value = '''3.30'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.64%  '

$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '''2.16'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +4.88%  '

$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").Value = '''2.75'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +4.08%  '

$ws.Range("D50").Value = '''2.91'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.12%  '

$ws.Range("D51").Value = '''144.32'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.06%  '
